$wb = $excel.ActiveWorkbook

$wsQuality = $wb.Worksheets.Item("quality_comparison")
$wsComp = $wb.Worksheets.Item("computational_comparison")

# --- quality_comparison sheet ---
# C1: top+bottom border (fontId 0) -> new style (borderId 4)
$c1 = $wsQuality.Range("C1")
$c1.ClearFormats()
$c1.Borders.LineStyle = 1
$c1.Borders.Item(7).LineStyle = -4142
$c1.Borders.Item(10).LineStyle = -4142

# D1: top+bottom+right border (fontId 0) -> new style (borderId 5)
$d1 = $wsQuality.Range("D1")
$d1.ClearFormats()
$d1.Borders.LineStyle = 1
$d1.Borders.Item(7).LineStyle = -4142

# C2 text: fedcore -> approach
$wsQuality.Range("C2").Value = "approach"

# --- computational_comparison sheet ---
# C1: top+bottom border
$c1b = $wsComp.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142

# D1: top+bottom+right border
$d1b = $wsComp.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

# F1: top+bottom border
$f1b = $wsComp.Range("F1")
$f1b.ClearFormats()
$f1b.Borders.LineStyle = 1
$f1b.Borders.Item(7).LineStyle = -4142
$f1b.Borders.Item(10).LineStyle = -4142

# G1: top+bottom+right border
$g1b = $wsComp.Range("G1")
$g1b.ClearFormats()
$g1b.Borders.LineStyle = 1
$g1b.Borders.Item(7).LineStyle = -4142

# C2 / F2 text: fedcore -> approach
$wsComp.Range("C2").Value = "approach"
$wsComp.Range("F2").Value = "approach"

# G5: remove the stray empty inline-string cell
$wsComp.Range("G5").ClearContents()
